# LoginPageTest: add test data rows for "Invalid entries in OTP field" script.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block (rows 14-18): OTP field header/value, and three more invalid-entry rows.
$ws.Range("A14").Value = "OTP field"
$ws.Range("B14").Value = 456789
$ws.Range("B15").Value = "abc123"
$ws.Range("B16").Value = "!@#$%^"
$ws.Range("B17").Value = 4567891
$ws.Range("B18").Value = 45678

# A14:B14 is the section header pair (mirrors A1/A7 styling) -- give it its own
# alignment-tagged style distinct from the "center" style already used by A1/A7.
$ws.Range("A14:B14").WrapText = $false

# B16 becomes a real hyperlink (mailto link), matching the existing B10 hyperlink.
$ws.Hyperlinks.Add($ws.Range("B16"), "mailto:!@#$%^")

# Move the active selection to the new last cell.
$null = $ws.Range("B18").Select()
